$d = $word.ActiveDocument

$d.Content.Find.Execute("208÷8=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "118÷6=19, 4", 2) | Out-Null
$d.Content.Find.Execute("244÷2=122, 0", $true, $false, $false, $false, $false, $true, 1, $false, "172÷8=21, 4", 2) | Out-Null
$d.Content.Find.Execute("897÷3=299, 0", $true, $false, $false, $false, $false, $true, 1, $false, "407÷2=203, 1", 2) | Out-Null
$d.Content.Find.Execute("420÷7=60, 0", $true, $false, $false, $false, $false, $true, 1, $false, "500÷3=166, 2", 2) | Out-Null
$d.Content.Find.Execute("502÷4=125, 2", $true, $false, $false, $false, $false, $true, 1, $false, "714÷7=102, 0", 2) | Out-Null
$d.Content.Find.Execute("530÷2=265, 0", $true, $false, $false, $false, $false, $true, 1, $false, "240÷8=30, 0", 2) | Out-Null
$d.Content.Find.Execute("805÷9=89, 4", $true, $false, $false, $false, $false, $true, 1, $false, "613÷3=204, 1", 2) | Out-Null
$d.Content.Find.Execute("830÷2=415, 0", $true, $false, $false, $false, $false, $true, 1, $false, "479÷9=53, 2", 2) | Out-Null
$d.Content.Find.Execute("341÷4=85, 1", $true, $false, $false, $false, $false, $true, 1, $false, "654÷5=130, 4", 2) | Out-Null
$d.Content.Find.Execute("769÷3=256, 1", $true, $false, $false, $false, $false, $true, 1, $false, "620÷8=77, 4", 2) | Out-Null
$d.Content.Find.Execute("270÷6=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "308÷4=77, 0", 2) | Out-Null
$d.Content.Find.Execute("205÷4=51, 1", $true, $false, $false, $false, $false, $true, 1, $false, "846÷7=120, 6", 2) | Out-Null
$d.Content.Find.Execute("458÷3=152, 2", $true, $false, $false, $false, $false, $true, 1, $false, "604÷2=302, 0", 2) | Out-Null
$d.Content.Find.Execute("293÷3=97, 2", $true, $false, $false, $false, $false, $true, 1, $false, "805÷7=115, 0", 2) | Out-Null
$d.Content.Find.Execute("525÷8=65, 5", $true, $false, $false, $false, $false, $true, 1, $false, "425÷4=106, 1", 2) | Out-Null
$d.Content.Find.Execute("628÷5=125, 3", $true, $false, $false, $false, $false, $true, 1, $false, "869÷3=289, 2", 2) | Out-Null
$d.Content.Find.Execute("679÷8=84, 7", $true, $false, $false, $false, $false, $true, 1, $false, "186÷3=62, 0", 2) | Out-Null
$d.Content.Find.Execute("187÷7=26, 5", $true, $false, $false, $false, $false, $true, 1, $false, "108÷7=15, 3", 2) | Out-Null
$d.Content.Find.Execute("502÷7=71, 5", $true, $false, $false, $false, $false, $true, 1, $false, "679÷9=75, 4", 2) | Out-Null
$d.Content.Find.Execute("712÷6=118, 4", $true, $false, $false, $false, $false, $true, 1, $false, "255÷5=51, 0", 2) | Out-Null
$d.Content.Find.Execute("105÷7=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "979÷9=108, 7", 2) | Out-Null
$d.Content.Find.Execute("173÷6=28, 5", $true, $false, $false, $false, $false, $true, 1, $false, "598÷3=199, 1", 2) | Out-Null
$d.Content.Find.Execute("665÷6=110, 5", $true, $false, $false, $false, $false, $true, 1, $false, "315÷9=35, 0", 2) | Out-Null
$d.Content.Find.Execute("665÷4=166, 1", $true, $false, $false, $false, $false, $true, 1, $false, "829÷3=276, 1", 2) | Out-Null
$d.Content.Find.Execute("154÷6=25, 4", $true, $false, $false, $false, $false, $true, 1, $false, "862÷8=107, 6", 2) | Out-Null
